# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the a2d5025a... row (row 4) on the zh-cn and de-de report sheets.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-03 10:14:55"   # Correspond Handoff Datetime
$wsZhCn.Range("G4").Value = "2016-03-03 10:15:49"   # Correspond Handback DateTime

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-03 10:15:08"   # Correspond Handoff Datetime
$wsDeDe.Range("G4").Value = "2016-03-03 10:16:14"   # Correspond Handback DateTime
